# ============================================================================
# KHL stats refresh: adds the 2025-12-16 match batch to Matches_SOG and rolls
# the derived shots-on-goal aggregates (Shots_HA / Shots_Summary) and Meta_ext
# forward to match. Mirrors the published "publish files + archive" commit.
# ============================================================================

$wb = $excel.ActiveWorkbook

# --- Matches_SOG: append the 8 new matches played 2025-12-16 (rows 373-380) ---
$wsMatches = $wb.Worksheets.Item("Matches_SOG")

$newMatches = @(
    @("897869", "2025-12-16T15:30:00", "Сибирь", "Авангард", 20, 37, "khl_text"),
    @("897868", "2025-12-16T17:00:00", "Автомобилист", "Барыс", 34, 27, "khl_text"),
    @("897871", "2025-12-16T18:00:00", "Лада", "ЦСКА", 30, 38, "khl_text"),
    @("897866", "2025-12-16T19:00:00", "Локомотив", "Динамо Мн", 25, 26, "khl_text"),
    @("897867", "2025-12-16T19:00:00", "Ак Барс", "Салават Юлаев", 27, 42, "khl_text"),
    @("897872", "2025-12-16T19:00:00", "Нефтехимик", "Торпедо", 19, 48, "khl_text"),
    @("897865", "2025-12-16T19:30:00", "Динамо М", "Спартак", 27, 26, "khl_text"),
    @("897870", "2025-12-16T19:30:00", "СКА", "Драконы", 36, 35, "khl_text")
)

$startRow = 373
for ($i = 0; $i -lt $newMatches.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $newMatches[$i]
    # column A (uid) is stored as text in this sheet, e.g. "897869" - force text
    # so the numeric-looking uid is not auto-coerced to a Number by COM.
    $wsMatches.Range("A$r").NumberFormat = "@"
    $wsMatches.Range("A$r").Value = $rowVals[0]
    $wsMatches.Range("B$r").Value = $rowVals[1]
    $wsMatches.Range("C$r").Value = $rowVals[2]
    $wsMatches.Range("D$r").Value = $rowVals[3]
    $wsMatches.Range("E$r").Value = $rowVals[4]
    $wsMatches.Range("F$r").Value = $rowVals[5]
    $wsMatches.Range("G$r").Value = $rowVals[6]
}

# --- Shots_HA: roll forward as_of_utc + home/away shots-on-goal aggregates ---
$wsShotsHA = $wb.Worksheets.Item("Shots_HA")

# row 2: D = "2025-12-16T19:30:00Z", E = 16, F = 16, G = 514, K = 568, L = 526, M = 35.5, N = 32.9
$wsShotsHA.Range("D2").Value = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("E2").Value = 16
$wsShotsHA.Range("F2").Value = 16
$wsShotsHA.Range("G2").Value = 514
$wsShotsHA.Range("K2").Value = 568
$wsShotsHA.Range("L2").Value = 526
$wsShotsHA.Range("M2").Value = 35.5
$wsShotsHA.Range("N2").Value = 32.9

# row 3: D = "2025-12-16T19:30:00Z", E = 15, F = 19, G = 436, H = 463, I = 29.1, J = 30.9, K = 536
$wsShotsHA.Range("D3").Value = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("E3").Value = 15
$wsShotsHA.Range("F3").Value = 19
$wsShotsHA.Range("G3").Value = 436
$wsShotsHA.Range("H3").Value = 463
$wsShotsHA.Range("I3").Value = 29.1
$wsShotsHA.Range("J3").Value = 30.9
$wsShotsHA.Range("K3").Value = 536

# row 4: D = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("D4").Value = "2025-12-16T19:30:00Z"

# row 5: D = "2025-12-16T19:30:00Z", E = 20, F = 17, G = 660, H = 554, I = 33, J = 27.7, K = 556
$wsShotsHA.Range("D5").Value = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("E5").Value = 20
$wsShotsHA.Range("F5").Value = 17
$wsShotsHA.Range("G5").Value = 660
$wsShotsHA.Range("H5").Value = 554
$wsShotsHA.Range("I5").Value = 33
$wsShotsHA.Range("J5").Value = 27.7
$wsShotsHA.Range("K5").Value = 556

# row 6: D = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("D6").Value = "2025-12-16T19:30:00Z"

# row 7: D = "2025-12-16T19:30:00Z", E = 21, F = 15, G = 682, K = 433, L = 492, M = 28.9, N = 32.8
$wsShotsHA.Range("D7").Value = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("E7").Value = 21
$wsShotsHA.Range("F7").Value = 15
$wsShotsHA.Range("G7").Value = 682
$wsShotsHA.Range("K7").Value = 433
$wsShotsHA.Range("L7").Value = 492
$wsShotsHA.Range("M7").Value = 28.9
$wsShotsHA.Range("N7").Value = 32.8

# row 8: D = "2025-12-16T19:30:00Z", E = 17, F = 17, G = 550, H = 452, I = 32.4, J = 26.6
$wsShotsHA.Range("D8").Value = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("E8").Value = 17
$wsShotsHA.Range("F8").Value = 17
$wsShotsHA.Range("G8").Value = 550
$wsShotsHA.Range("H8").Value = 452
$wsShotsHA.Range("I8").Value = 32.4
$wsShotsHA.Range("J8").Value = 26.6

# row 9: D = "2025-12-16T19:30:00Z", E = 18, F = 15, G = 661, K = 519, L = 418, M = 34.6, N = 27.9
$wsShotsHA.Range("D9").Value = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("E9").Value = 18
$wsShotsHA.Range("F9").Value = 15
$wsShotsHA.Range("G9").Value = 661
$wsShotsHA.Range("K9").Value = 519
$wsShotsHA.Range("L9").Value = 418
$wsShotsHA.Range("M9").Value = 34.6
$wsShotsHA.Range("N9").Value = 27.9

# row 10: D = "2025-12-16T19:30:00Z", E = 14, F = 20, G = 408, K = 566, L = 725, M = 28.3, N = 36.2
$wsShotsHA.Range("D10").Value = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("E10").Value = 14
$wsShotsHA.Range("F10").Value = 20
$wsShotsHA.Range("G10").Value = 408
$wsShotsHA.Range("K10").Value = 566
$wsShotsHA.Range("L10").Value = 725
$wsShotsHA.Range("M10").Value = 28.3
$wsShotsHA.Range("N10").Value = 36.2

# row 11: D = "2025-12-16T19:30:00Z", E = 16, F = 18, G = 422, H = 574, I = 26.4, J = 35.9, K = 435
$wsShotsHA.Range("D11").Value = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("E11").Value = 16
$wsShotsHA.Range("F11").Value = 18
$wsShotsHA.Range("G11").Value = 422
$wsShotsHA.Range("H11").Value = 574
$wsShotsHA.Range("I11").Value = 26.4
$wsShotsHA.Range("J11").Value = 35.9
$wsShotsHA.Range("K11").Value = 435

# row 12: D = "2025-12-16T19:30:00Z", E = 18, F = 19, G = 597, H = 476, I = 33.2, J = 26.4, K = 576
$wsShotsHA.Range("D12").Value = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("E12").Value = 18
$wsShotsHA.Range("F12").Value = 19
$wsShotsHA.Range("G12").Value = 597
$wsShotsHA.Range("H12").Value = 476
$wsShotsHA.Range("I12").Value = 33.2
$wsShotsHA.Range("J12").Value = 26.4
$wsShotsHA.Range("K12").Value = 576

# row 13: D = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("D13").Value = "2025-12-16T19:30:00Z"

# row 14: D = "2025-12-16T19:30:00Z", E = 22, F = 13, G = 674, H = 785, I = 30.6, J = 35.7, K = 374
$wsShotsHA.Range("D14").Value = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("E14").Value = 22
$wsShotsHA.Range("F14").Value = 13
$wsShotsHA.Range("G14").Value = 674
$wsShotsHA.Range("H14").Value = 785
$wsShotsHA.Range("I14").Value = 30.6
$wsShotsHA.Range("J14").Value = 35.7
$wsShotsHA.Range("K14").Value = 374

# row 15: D = "2025-12-16T19:30:00Z", E = 17, F = 17, G = 554, H = 578, I = 32.6, J = 34, K = 513
$wsShotsHA.Range("D15").Value = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("E15").Value = 17
$wsShotsHA.Range("F15").Value = 17
$wsShotsHA.Range("G15").Value = 554
$wsShotsHA.Range("H15").Value = 578
$wsShotsHA.Range("I15").Value = 32.6
$wsShotsHA.Range("J15").Value = 34
$wsShotsHA.Range("K15").Value = 513

# row 16: D = "2025-12-16T19:30:00Z", E = 14, F = 21, G = 389, K = 614, L = 645, M = 29.2, N = 30.7
$wsShotsHA.Range("D16").Value = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("E16").Value = 14
$wsShotsHA.Range("F16").Value = 21
$wsShotsHA.Range("G16").Value = 389
$wsShotsHA.Range("K16").Value = 614
$wsShotsHA.Range("L16").Value = 645
$wsShotsHA.Range("M16").Value = 29.2
$wsShotsHA.Range("N16").Value = 30.7

# row 17: D = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("D17").Value = "2025-12-16T19:30:00Z"

# row 18: D = "2025-12-16T19:30:00Z", E = 16, F = 19, G = 423, H = 545, I = 26.4, J = 34.1, K = 500
$wsShotsHA.Range("D18").Value = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("E18").Value = 16
$wsShotsHA.Range("F18").Value = 19
$wsShotsHA.Range("G18").Value = 423
$wsShotsHA.Range("H18").Value = 545
$wsShotsHA.Range("I18").Value = 26.4
$wsShotsHA.Range("J18").Value = 34.1
$wsShotsHA.Range("K18").Value = 500

# row 19: D = "2025-12-16T19:30:00Z", E = 19, F = 16, G = 641, K = 508, L = 526, M = 31.8, N = 32.9
$wsShotsHA.Range("D19").Value = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("E19").Value = 19
$wsShotsHA.Range("F19").Value = 16
$wsShotsHA.Range("G19").Value = 641
$wsShotsHA.Range("K19").Value = 508
$wsShotsHA.Range("L19").Value = 526
$wsShotsHA.Range("M19").Value = 31.8
$wsShotsHA.Range("N19").Value = 32.9

# row 20: D = "2025-12-16T19:30:00Z", E = 17, F = 19, G = 565, K = 645, L = 577, M = 33.9, N = 30.4
$wsShotsHA.Range("D20").Value = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("E20").Value = 17
$wsShotsHA.Range("F20").Value = 19
$wsShotsHA.Range("G20").Value = 565
$wsShotsHA.Range("K20").Value = 645
$wsShotsHA.Range("L20").Value = 577
$wsShotsHA.Range("M20").Value = 33.9
$wsShotsHA.Range("N20").Value = 30.4

# row 21: D = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("D21").Value = "2025-12-16T19:30:00Z"

# row 22: D = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("D22").Value = "2025-12-16T19:30:00Z"

# row 23: D = "2025-12-16T19:30:00Z", E = 16, F = 19, G = 426, K = 494, L = 504, M = 26, N = 26.5
$wsShotsHA.Range("D23").Value = "2025-12-16T19:30:00Z"
$wsShotsHA.Range("E23").Value = 16
$wsShotsHA.Range("F23").Value = 19
$wsShotsHA.Range("G23").Value = 426
$wsShotsHA.Range("K23").Value = 494
$wsShotsHA.Range("L23").Value = 504
$wsShotsHA.Range("M23").Value = 26
$wsShotsHA.Range("N23").Value = 26.5

# --- Shots_Summary: roll forward as_of_utc + total shots-on-goal aggregates ---
$wsShotsSummary = $wb.Worksheets.Item("Shots_Summary")

# row 2: D = "2025-12-16T19:30:00Z", E = 32, F = 1082, G = 987, H = 33.8, I = 30.8
$wsShotsSummary.Range("D2").Value = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("E2").Value = 32
$wsShotsSummary.Range("F2").Value = 1082
$wsShotsSummary.Range("G2").Value = 987
$wsShotsSummary.Range("H2").Value = 33.8
$wsShotsSummary.Range("I2").Value = 30.8

# row 3: D = "2025-12-16T19:30:00Z", E = 34, F = 972, G = 1062, H = 28.6, I = 31.2
$wsShotsSummary.Range("D3").Value = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("E3").Value = 34
$wsShotsSummary.Range("F3").Value = 972
$wsShotsSummary.Range("G3").Value = 1062
$wsShotsSummary.Range("H3").Value = 28.6
$wsShotsSummary.Range("I3").Value = 31.2

# row 4: D = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("D4").Value = "2025-12-16T19:30:00Z"

# row 5: D = "2025-12-16T19:30:00Z", E = 37, F = 1216, G = 1048, H = 32.9, I = 28.3
$wsShotsSummary.Range("D5").Value = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("E5").Value = 37
$wsShotsSummary.Range("F5").Value = 1216
$wsShotsSummary.Range("G5").Value = 1048
$wsShotsSummary.Range("H5").Value = 32.9
$wsShotsSummary.Range("I5").Value = 28.3

# row 6: D = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("D6").Value = "2025-12-16T19:30:00Z"

# row 7: D = "2025-12-16T19:30:00Z", E = 36, F = 1115, G = 1140, H = 31, I = 31.7
$wsShotsSummary.Range("D7").Value = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("E7").Value = 36
$wsShotsSummary.Range("F7").Value = 1115
$wsShotsSummary.Range("G7").Value = 1140
$wsShotsSummary.Range("H7").Value = 31
$wsShotsSummary.Range("I7").Value = 31.7

# row 8: D = "2025-12-16T19:30:00Z", E = 34, F = 1046, G = 966, H = 30.8, I = 28.4
$wsShotsSummary.Range("D8").Value = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("E8").Value = 34
$wsShotsSummary.Range("F8").Value = 1046
$wsShotsSummary.Range("G8").Value = 966
$wsShotsSummary.Range("H8").Value = 30.8
$wsShotsSummary.Range("I8").Value = 28.4

# row 9: D = "2025-12-16T19:30:00Z", E = 33, F = 1180, G = 901, H = 35.8, I = 27.3
$wsShotsSummary.Range("D9").Value = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("E9").Value = 33
$wsShotsSummary.Range("F9").Value = 1180
$wsShotsSummary.Range("G9").Value = 901
$wsShotsSummary.Range("H9").Value = 35.8
$wsShotsSummary.Range("I9").Value = 27.3

# row 10: D = "2025-12-16T19:30:00Z", E = 34, F = 974, G = 1213, H = 28.6, I = 35.7
$wsShotsSummary.Range("D10").Value = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("E10").Value = 34
$wsShotsSummary.Range("F10").Value = 974
$wsShotsSummary.Range("G10").Value = 1213
$wsShotsSummary.Range("H10").Value = 28.6
$wsShotsSummary.Range("I10").Value = 35.7

# row 11: D = "2025-12-16T19:30:00Z", E = 34, F = 857, G = 1257, H = 25.2, I = 37
$wsShotsSummary.Range("D11").Value = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("E11").Value = 34
$wsShotsSummary.Range("F11").Value = 857
$wsShotsSummary.Range("G11").Value = 1257
$wsShotsSummary.Range("H11").Value = 25.2
$wsShotsSummary.Range("I11").Value = 37

# row 12: D = "2025-12-16T19:30:00Z", E = 37, F = 1173, G = 946, H = 31.7, I = 25.6
$wsShotsSummary.Range("D12").Value = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("E12").Value = 37
$wsShotsSummary.Range("F12").Value = 1173
$wsShotsSummary.Range("G12").Value = 946
$wsShotsSummary.Range("H12").Value = 31.7
$wsShotsSummary.Range("I12").Value = 25.6

# row 13: D = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("D13").Value = "2025-12-16T19:30:00Z"

# row 14: D = "2025-12-16T19:30:00Z", E = 35, F = 1048, G = 1260, H = 29.9, I = 36
$wsShotsSummary.Range("D14").Value = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("E14").Value = 35
$wsShotsSummary.Range("F14").Value = 1048
$wsShotsSummary.Range("G14").Value = 1260
$wsShotsSummary.Range("H14").Value = 29.9
$wsShotsSummary.Range("I14").Value = 36

# row 15: D = "2025-12-16T19:30:00Z", E = 34, F = 1067, G = 1161, H = 31.4, I = 34.1
$wsShotsSummary.Range("D15").Value = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("E15").Value = 34
$wsShotsSummary.Range("F15").Value = 1067
$wsShotsSummary.Range("G15").Value = 1161
$wsShotsSummary.Range("H15").Value = 31.4
$wsShotsSummary.Range("I15").Value = 34.1

# row 16: D = "2025-12-16T19:30:00Z", E = 35, F = 1003, G = 1034, H = 28.7, I = 29.5
$wsShotsSummary.Range("D16").Value = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("E16").Value = 35
$wsShotsSummary.Range("F16").Value = 1003
$wsShotsSummary.Range("G16").Value = 1034
$wsShotsSummary.Range("H16").Value = 28.7
$wsShotsSummary.Range("I16").Value = 29.5

# row 17: D = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("D17").Value = "2025-12-16T19:30:00Z"

# row 18: D = "2025-12-16T19:30:00Z", E = 35, F = 923, G = 1190, H = 26.4, I = 34
$wsShotsSummary.Range("D18").Value = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("E18").Value = 35
$wsShotsSummary.Range("F18").Value = 923
$wsShotsSummary.Range("G18").Value = 1190
$wsShotsSummary.Range("H18").Value = 26.4
$wsShotsSummary.Range("I18").Value = 34

# row 19: D = "2025-12-16T19:30:00Z", E = 35, F = 1149, G = 1068, H = 32.8, I = 30.5
$wsShotsSummary.Range("D19").Value = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("E19").Value = 35
$wsShotsSummary.Range("F19").Value = 1149
$wsShotsSummary.Range("G19").Value = 1068
$wsShotsSummary.Range("H19").Value = 32.8
$wsShotsSummary.Range("I19").Value = 30.5

# row 20: D = "2025-12-16T19:30:00Z", E = 36, F = 1210, G = 1107, H = 33.6, I = 30.8
$wsShotsSummary.Range("D20").Value = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("E20").Value = 36
$wsShotsSummary.Range("F20").Value = 1210
$wsShotsSummary.Range("G20").Value = 1107
$wsShotsSummary.Range("H20").Value = 33.6
$wsShotsSummary.Range("I20").Value = 30.8

# row 21: D = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("D21").Value = "2025-12-16T19:30:00Z"

# row 22: D = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("D22").Value = "2025-12-16T19:30:00Z"

# row 23: D = "2025-12-16T19:30:00Z", E = 35, F = 920, G = 939, H = 26.3, I = 26.8
$wsShotsSummary.Range("D23").Value = "2025-12-16T19:30:00Z"
$wsShotsSummary.Range("E23").Value = 35
$wsShotsSummary.Range("F23").Value = 920
$wsShotsSummary.Range("G23").Value = 939
$wsShotsSummary.Range("H23").Value = 26.3
$wsShotsSummary.Range("I23").Value = 26.8

# --- Meta_ext: bump as_of_utc + build_version ---
$wsMeta = $wb.Worksheets.Item("Meta_ext")
$wsMeta.Range("B2").Value = "2025-12-16T19:30:00Z"
$wsMeta.Range("D2").Value = 59
